{"js": "// Author's commit: merges the split date-range runs (\"2019 - \" + \"present\"\n// and \"2019 - \" + \"2020\", each originally followed by a third run holding\n// the tab + following label text) into one run apiece, with no visible text\n// change - and updates the Scopus citation count from 590 to 615 (h-index\n// stays at 13).\n//\n// Body.search() matches text across run boundaries; Range.insertText(text,\n// \"Replace\") rewrites the matched range in place (re-merging it into a\n// single run). The replacement text is identical to the current text, so\n// the only effect is the run merge itself.\n\nconst body = context.document.body;\n\n// 1) \"2019 - \" / \"present\" / \"\\tPI, ESR Strategic Science Investment\n//    Funding: \" runs -> a single run.\nconst presentText =\n  \"2019 - present\\tPI, ESR Strategic Science Investment Funding: \";\nconst presentResults = body.search(presentText, { matchCase: true });\npresentResults.load(\"items\");\n\n// 2) \"2019 - \" / \"2020\" / \"\\tPI, ESR Pioneer Funding: \" runs -> a single run.\nconst twentyTwentyText = \"2019 - 2020\\tPI, ESR Pioneer Funding: \";\nconst twentyTwentyResults = body.search(twentyTwentyText, {\n  matchCase: true,\n});\ntwentyTwentyResults.load(\"items\");\n\nawait context.sync();\n\nif (presentResults.items.length > 0) {\n  presentResults.items[0].insertText(presentText, \"Replace\");\n}\nif (twentyTwentyResults.items.length > 0) {\n  twentyTwentyResults.items[0].insertText(twentyTwentyText, \"Replace\");\n}\nawait context.sync();\n\n// 3) Scopus citation count: \"h-index of 13 (cited 590 times)\" -> \"... 615 times)\".\n//    Split the edit in two so the \"615\" digits stay in their own run and keep\n//    their distinct (en-NZ) run formatting, matching the source edit exactly:\n//      - trim the trailing \"5\" off \"...(cited 5\"\n//      - change the adjoining \"90\" run to \"615\"\nconst citedResults = body.search(\"(cited 5\", { matchCase: true });\ncitedResults.load(\"items\");\nawait context.sync();\nif (citedResults.items.length > 0) {\n  citedResults.items[0].insertText(\"(cited \", \"Replace\");\n}\nawait context.sync();\n\nconst digitsResults = body.search(\"90\", { matchCase: true });\ndigitsResults.load(\"items\");\nawait context.sync();\nif (digitsResults.items.length > 0) {\n  digitsResults.items[0].insertText(\"615\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Author's commit: merges the split date-range runs (\"2019 - \" + \"present\"\n# and \"2019 - \" + \"2020\") into single runs with no visible text change, and\n# updates the Scopus citation count from 590 to 615 (h-index stays at 13).\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n# 1) \"2019 - \" / \"present\" runs -> single run \"2019 - present\"\n$find1 = $d.Content.Find\n$find1.Execute(\"2019 - present\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"2019 - present\", $wdReplaceOne) | Out-Null\n\n# 2) \"2019 - \" / \"2020\" runs -> single run \"2019 - 2020\"\n$find2 = $d.Content.Find\n$find2.Execute(\"2019 - 2020\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"2019 - 2020\", $wdReplaceOne) | Out-Null\n\n# 3) Scopus citation count: \"h-index of 13 (cited 590 times)\" -> \"... 615 times)\".\n#    Keep the digits run (\"90\" -> \"615\") separate from the surrounding text so the\n#    distinct run-level formatting on the number is preserved, matching the source edit.\n$find3 = $d.Content.Find\n$find3.Execute(\"(cited 5\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"(cited \", $wdReplaceOne) | Out-Null\n\n$find4 = $d.Content.Find\n$find4.Execute(\"90\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"615\", $wdReplaceOne) | Out-Null\n"}
